# Generate Report for Handoff
# Replaces the two tracked files (898253a8... / b4576177...) with a new
# pair (8cf671d6... / ffff8ac5cfcc...) that is "Ready for handoff" instead
# of "Handed back: in sync with en-US", and refreshes all the associated
# dates, target-file hyperlinks, and per-language detail rows.

$wb = $excel.ActiveWorkbook

$newFile1 = "8cf671d6-5884-4d8e-9ae8-d2acbb21284d.md"
$newFile2 = "ffff8ac5cfcc-76db-407d-81e8-5e24277676cc.md"
$newHash  = "d64bffcf65a9b190c857baf1334016a6b9a9554f"
$newUuid1 = "8cf671d6-5884-4d8e-9ae8-d2acbb21284d"

$zhXlf = "$newUuid1.$newHash.zh-cn.xlf"
$deXlf = "$newUuid1.$newHash.de-de.xlf"

$status = "Ready for handoff"
$overviewDate = "2016-06-13 17:06:51"
$zhHandoffDatetime = "2016-03-13 17:06:47"
$deHandoffDatetime = "2016-03-13 17:06:51"
$handbackDatetime = "0001-01-01 00:00:00"
$include = "Include"
$mdExt = ".md"

$srcBase = "https://github.com/OpenLocalizationTest/oltest/blob/e3b12d17401494a259b2b0bbd15128741d23416b/e2e"
$zhOffBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c2f953043dad3a1ab1316429d32f1f91fbfab23d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht"
$deOffBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8cef05287ee3f17324125c53109364eca9e17010/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht"

# ---------------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = $newFile1
$ws.Range("B2").Value = $status
$ws.Range("C2").Value = $status
$ws.Range("D2").Value = $overviewDate

$ws.Range("A3").Value = $newFile2
$ws.Range("B3").Value = $status
$ws.Range("C3").Value = $status
$ws.Range("D3").Value = $overviewDate

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "$srcBase/$newFile1", "", "", $newFile1)
$ws.Hyperlinks.Add($ws.Range("A3"), "$srcBase/$newFile2", "", "", $newFile2)

# ---------------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value = $newFile1
$ws.Range("B2").Value = $mdExt
$ws.Range("C2").Value = $status
$ws.Range("D2").Value = $zhXlf
$ws.Range("E2").Value = $zhHandoffDatetime
$ws.Range("H2").Value = $handbackDatetime
$ws.Range("I2").Value = $include

$ws.Range("A3").Value = $newFile2
$ws.Range("B3").Value = $mdExt
$ws.Range("C3").Value = $status
$ws.Range("D3").Value = $zhXlf
$ws.Range("E3").Value = $zhHandoffDatetime
$ws.Range("H3").Value = $handbackDatetime
$ws.Range("I3").Value = $include

# Drop the now-unused "Latest Target File" / "Latest Handback File" columns
# for both data rows.
$ws.Range("F2:G3").Clear()

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "$srcBase/$newFile1", "", "", $newFile1)
$ws.Hyperlinks.Add($ws.Range("B2"), "$srcBase/$newFile1", "", "", $mdExt)
$ws.Hyperlinks.Add($ws.Range("D2"), "$zhOffBase/$zhXlf", "", "", $zhXlf)
$ws.Hyperlinks.Add($ws.Range("A3"), "$srcBase/$newFile2", "", "", $newFile2)
$ws.Hyperlinks.Add($ws.Range("B3"), "$srcBase/$newFile2", "", "", $mdExt)
$ws.Hyperlinks.Add($ws.Range("D3"), "$zhOffBase/$zhXlf", "", "", $zhXlf)

# ---------------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value = $newFile1
$ws.Range("B2").Value = $mdExt
$ws.Range("C2").Value = $status
$ws.Range("D2").Value = $deXlf
$ws.Range("E2").Value = $deHandoffDatetime
$ws.Range("H2").Value = $handbackDatetime
$ws.Range("I2").Value = $include

$ws.Range("A3").Value = $newFile2
$ws.Range("B3").Value = $mdExt
$ws.Range("C3").Value = $status
$ws.Range("D3").Value = $deXlf
$ws.Range("E3").Value = $deHandoffDatetime
$ws.Range("H3").Value = $handbackDatetime
$ws.Range("I3").Value = $include

$ws.Range("F2:G3").Clear()

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "$srcBase/$newFile1", "", "", $newFile1)
$ws.Hyperlinks.Add($ws.Range("B2"), "$srcBase/$newFile1", "", "", $mdExt)
$ws.Hyperlinks.Add($ws.Range("D2"), "$deOffBase/$deXlf", "", "", $deXlf)
$ws.Hyperlinks.Add($ws.Range("A3"), "$srcBase/$newFile2", "", "", $newFile2)
$ws.Hyperlinks.Add($ws.Range("B3"), "$srcBase/$newFile2", "", "", $mdExt)
$ws.Hyperlinks.Add($ws.Range("D3"), "$deOffBase/$deXlf", "", "", $deXlf)
